$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update A18 content to new URL and add a hyperlink to it
$ws.Range("A18").Value = "http://www.southcn.com/"
$ws.Hyperlinks.Add($ws.Range("A18"), "http://www.southcn.com/")

# A19 and A20 now carry the "yellow fill" formatting used elsewhere in the sheet
$ws.Range("A19").Interior.Color = 65535
$ws.Range("A20").Interior.Color = 65535

# Update the selection / scrolled position of the sheet view
$ws.Range("A13").Select()
$excel.ActiveWindow.ScrollRow = 13
$ws.Range("E27").Select()
